# Update the "想去人数" (interest count) values in column F across the
# four sheets of the 广州-漫展信息 workbook, per the generated-output refresh.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# 展览
$ws1.Range("F3").Value = 244
$ws1.Range("F5").Value = 241
$ws1.Range("F6").Value = 412
$ws1.Range("F7").Value = 588
$ws1.Range("F10").Value = 342
$ws1.Range("F11").Value = 145
$ws1.Range("F12").Value = 651
$ws1.Range("F13").Value = 86
$ws1.Range("F14").Value = 1801
$ws1.Range("F15").Value = 348
$ws1.Range("F16").Value = 3024
$ws1.Range("F17").Value = 316
$ws1.Range("F18").Value = 493
$ws1.Range("F20").Value = 142

# 演出
$ws2.Range("F5").Value = 19
$ws2.Range("F13").Value = 91
$ws2.Range("F14").Value = 40

# 本地生活
$ws3.Range("F2").Value = 5318
$ws3.Range("F3").Value = 318
$ws3.Range("F4").Value = 246

# 全部类型
$ws4.Range("F3").Value = 5318
$ws4.Range("F4").Value = 318
$ws4.Range("F6").Value = 246
$ws4.Range("F7").Value = 244
$ws4.Range("F10").Value = 19
$ws4.Range("F16").Value = 241
$ws4.Range("F17").Value = 412
$ws4.Range("F18").Value = 588
$ws4.Range("F22").Value = 342
$ws4.Range("F23").Value = 145
$ws4.Range("F26").Value = 651
$ws4.Range("F27").Value = 86
$ws4.Range("F28").Value = 91
$ws4.Range("F29").Value = 1801
$ws4.Range("F30").Value = 348
$ws4.Range("F31").Value = 3025
$ws4.Range("F32").Value = 40
$ws4.Range("F33").Value = 316
$ws4.Range("F34").Value = 493
$ws4.Range("F37").Value = 142
